# Update Excel SCD0011 until SCD0016
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0185 to SCD0011
$ws.Name = "SCD0011"

# Update TC_ID cell (B2) from "DGS-200" to "SCD0011-016"
$ws.Range("B2").Value = "SCD0011-016"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns("B:B").ColumnWidth = 11.6

# Move the active selection to B3
$ws.Range("B3").Select() | Out-Null
